$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046015360399743
$ws.Range("D2").Value = 1.045517149817485
$ws.Range("E2").Value = 1.057916378128792
$ws.Range("F2").Value = 1.065189534509206
$ws.Range("I2").Value = 1.041297323242586
$ws.Range("J2").Value = 1.051072011530409
$ws.Range("K2").Value = 1.048285168081746
$ws.Range("L2").Value = 1.060650038855828
$ws.Range("M2").Value = 1.067903444620381
$ws.Range("N2").Value = 1.052564654346427

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047145190270167
$ws.Range("D3").Value = 1.04635108905697
$ws.Range("E3").Value = 1.059277940353205
$ws.Range("F3").Value = 1.066647310122513
$ws.Range("I3").Value = 1.041598539744753
$ws.Range("J3").Value = 1.051848848591435
$ws.Range("K3").Value = 1.048930505235321
$ws.Range("L3").Value = 1.061824142150128
$ws.Range("M3").Value = 1.069174965131756
$ws.Range("N3").Value = 1.053342594605184

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047875695047077
$ws.Range("D4").Value = 1.04689020647245
$ws.Range("E4").Value = 1.060159221471009
$ws.Range("F4").Value = 1.067590910002192
$ws.Range("I4").Value = 1.04179200307733
$ws.Range("J4").Value = 1.052350391061534
$ws.Range("K4").Value = 1.049346939363912
$ws.Range("L4").Value = 1.062583587038346
$ws.Range("M4").Value = 1.069997524829581
$ws.Range("N4").Value = 1.053844849323119

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048182664422545
$ws.Range("D5").Value = 1.047116733025788
$ws.Range("E5").Value = 1.060529777421038
$ws.Range("F5").Value = 1.067987679900764
$ws.Range("I5").Value = 1.041872989818786
$ws.Range("J5").Value = 1.052560971752177
$ws.Range("K5").Value = 1.049521735587783
$ws.Range("L5").Value = 1.062902793931102
$ws.Range("M5").Value = 1.070343284330038
$ws.Range("N5").Value = 1.054055729062496

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048234198076829
$ws.Range("D6").Value = 1.047154760904196
$ws.Range("E6").Value = 1.060591999324335
$ws.Range("F6").Value = 1.068054304136062
$ws.Range("I6").Value = 1.041886567616127
$ws.Range("J6").Value = 1.052596313506763
$ws.Range("K6").Value = 1.049551068674385
$ws.Range("L6").Value = 1.062956386520299
$ws.Range("M6").Value = 1.07040133634556
$ws.Range("N6").Value = 1.054091121006427

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047879797314457
$ws.Range("D7").Value = 1.046893233795888
$ws.Range("E7").Value = 1.060164172596788
$ws.Range("F7").Value = 1.067596211340247
$ws.Range("I7").Value = 1.041793086582176
$ws.Range("J7").Value = 1.052353205901795
$ws.Range("K7").Value = 1.049349276070142
$ws.Range("L7").Value = 1.062587852543267
$ws.Range("M7").Value = 1.070002145056352
$ws.Range("N7").Value = 1.053847668160776

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046397310942408
$ws.Range("D8").Value = 1.045799086756008
$ws.Range("E8").Value = 1.058376472744236
$ws.Range("F8").Value = 1.065682132138577
$ws.Range("I8").Value = 1.041399419992913
$ws.Range("J8").Value = 1.051334780450333
$ws.Range("K8").Value = 1.048503499740848
$ws.Range("L8").Value = 1.061046891977262
$ws.Range("M8").Value = 1.068333204001855
$ws.Range("N8").Value = 1.052827796428358

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043780528876924
$ws.Range("D9").Value = 1.043867226315437
$ws.Range("E9").Value = 1.055228156001323
$ws.Range("F9").Value = 1.062311585217766
$ws.Range("I9").Value = 1.040694651362952
$ws.Range("J9").Value = 1.049531533768213
$ws.Range("K9").Value = 1.047004351974475
$ws.Range("L9").Value = 1.058329259004191
$ws.Range("M9").Value = 1.065390654363388
$ws.Range("N9").Value = 1.051021988929114

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042032898884194
$ws.Range("D10").Value = 1.042576701458593
$ws.Range("E10").Value = 1.053130307153263
$ws.Range("F10").Value = 1.060065881721935
$ws.Range("I10").Value = 1.040217329108679
$ws.Range("J10").Value = 1.048323481143437
$ws.Range("K10").Value = 1.045998965655584
$ws.Range("L10").Value = 1.056515797289398
$ws.Range("M10").Value = 1.063427643272583
$ws.Range("N10").Value = 1.049812220731045

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041275391949081
$ws.Range("D11").Value = 1.042017259125627
$ws.Range("E11").Value = 1.0522221129009
$ws.Range("F11").Value = 1.059093732334665
$ws.Range("I11").Value = 1.040008863219521
$ws.Range("J11").Value = 1.047798967711719
$ws.Range("K11").Value = 1.04556219810629
$ws.Range("L11").Value = 1.055730102108138
$ws.Range("M11").Value = 1.06257728315953
$ws.Range("N11").Value = 1.049286962430091

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040993901827502
$ws.Range("D12").Value = 1.041809360742421
$ws.Range("E12").Value = 1.05188479330139
$ws.Range("F12").Value = 1.058732666712874
$ws.Range("I12").Value = 1.039931161302797
$ws.Range("J12").Value = 1.047603925391607
$ws.Range("K12").Value = 1.045399747440971
$ws.Range("L12").Value = 1.055438188037018
$ws.Range("M12").Value = 1.062261363128128
$ws.Range("N12").Value = 1.049091643127513

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041054287806758
$ws.Range("D13").Value = 1.041853960037779
$ws.Range("E13").Value = 1.05195714847189
$ws.Range("F13").Value = 1.058810115048641
$ws.Range("I13").Value = 1.03994784078396
$ws.Range("J13").Value = 1.047645772370754
$ws.Range("K13").Value = 1.045434603441255
$ws.Range("L13").Value = 1.055500807908364
$ws.Range("M13").Value = 1.062329131752709
$ws.Range("N13").Value = 1.04913354953417

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041252126315333
$ws.Range("D14").Value = 1.042000076156422
$ws.Range("E14").Value = 1.052194229487904
$ws.Range("F14").Value = 1.059063885870584
$ws.Range("I14").Value = 1.040002445833954
$ws.Range("J14").Value = 1.047782849840976
$ws.Range("K14").Value = 1.045548774291328
$ws.Range("L14").Value = 1.055705973887768
$ws.Range("M14").Value = 1.062551170325203
$ws.Range("N14").Value = 1.049270821670124

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0413740054
$ws.Range("D15").Value = 1.042090090291463
$ws.Range("E15").Value = 1.052340306008705
$ws.Range("F15").Value = 1.059220246816607
$ws.Range("I15").Value = 1.040036054221947
$ws.Range("J15").Value = 1.047867279301499
$ws.Range("K15").Value = 1.04561909009771
$ws.Range("L15").Value = 1.055832373782651
$ws.Range("M15").Value = 1.062687967761451
$ws.Range("N15").Value = 1.049355371030164

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042083155663737
$ws.Range("D16").Value = 1.042613816310864
$ws.Range("E16").Value = 1.053190584630412
$ws.Range("F16").Value = 1.060130405075681
$ws.Range("I16").Value = 1.040231126698994
$ws.Range("J16").Value = 1.048358261351549
$ws.Range("K16").Value = 1.046027922284708
$ws.Range("L16").Value = 1.056567931408603
$ws.Range("M16").Value = 1.063484070890263
$ws.Range("N16").Value = 1.049847050331042

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042527778620824
$ws.Range("D17").Value = 1.042942164710765
$ws.Range("E17").Value = 1.053723989438937
$ws.Range("F17").Value = 1.060701388506936
$ws.Range("I17").Value = 1.040353012907737
$ws.Range("J17").Value = 1.048665860518172
$ws.Range("K17").Value = 1.046283988643821
$ws.Range("L17").Value = 1.057029203077814
$ws.Range("M17").Value = 1.063983345367746
$ws.Range("N17").Value = 1.050155086323763

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042787045380509
$ws.Range("D18").Value = 1.043133623423709
$ws.Range("E18").Value = 1.054035133866791
$ws.Range("F18").Value = 1.061034458384277
$ws.Range("I18").Value = 1.040423935158002
$ws.Range("J18").Value = 1.048845141130115
$ws.Range("K18").Value = 1.046433210050985
$ws.Range("L18").Value = 1.057298211450857
$ws.Range("M18").Value = 1.06427452892241
$ws.Range("N18").Value = 1.050334621534739

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042875436117435
$ws.Range("D19").Value = 1.043198895519417
$ws.Range("E19").Value = 1.054141229325654
$ws.Range("F19").Value = 1.0611480309618
$ws.Range("I19").Value = 1.040448088672748
$ws.Range("J19").Value = 1.048906248032538
$ws.Range("K19").Value = 1.046484067360554
$ws.Range("L19").Value = 1.057389929071813
$ws.Range("M19").Value = 1.064373809203272
$ws.Range("N19").Value = 1.050395815215972

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042480082486382
$ws.Range("D20").Value = 1.042906942387696
$ws.Range("E20").Value = 1.053666758274075
$ws.Range("F20").Value = 1.06064012484385
$ws.Range("I20").Value = 1.040339953457406
$ws.Range("J20").Value = 1.048632872197119
$ws.Range("K20").Value = 1.046256529404693
$ws.Range("L20").Value = 1.056979717527527
$ws.Range("M20").Value = 1.063929781568285
$ws.Range("N20").Value = 1.05012205115551

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041193871090958
$ws.Range("D21").Value = 1.041957051275459
$ws.Range("E21").Value = 1.052124414403052
$ws.Range("F21").Value = 1.058989155786703
$ws.Range("I21").Value = 1.039986373422195
$ws.Range("J21").Value = 1.047742489881107
$ws.Range("K21").Value = 1.045515159786363
$ws.Range("L21").Value = 1.055645559634516
$ws.Range("M21").Value = 1.062485787143369
$ws.Range("N21").Value = 1.049230404394481

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040384493195353
$ws.Range("D22").Value = 1.041359257702614
$ws.Range("E22").Value = 1.051154819482308
$ws.Range("F22").Value = 1.057951320238755
$ws.Range("I22").Value = 1.039762510397358
$ws.Range("J22").Value = 1.047181428354082
$ws.Range("K22").Value = 1.045047782460302
$ws.Range("L22").Value = 1.054806303824465
$ws.Range("M22").Value = 1.061577550612767
$ws.Range("N22").Value = 1.048668546095733

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040813625533458
$ws.Range("D23").Value = 1.041676212726236
$ws.Range("E23").Value = 1.051668808391941
$ws.Range("F23").Value = 1.058501479552551
$ws.Range("I23").Value = 1.039881331867267
$ws.Range("J23").Value = 1.047478976006601
$ws.Range("K23").Value = 1.045295666771765
$ws.Range("L23").Value = 1.055251250030561
$ws.Range("M23").Value = 1.062059057485951
$ws.Range("N23").Value = 1.048966516300048

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.04250163454539
$ws.Range("D24").Value = 1.042922858028608
$ws.Range("E24").Value = 1.053692618516088
$ws.Range("F24").Value = 1.06066780717816
$ws.Range("I24").Value = 1.040345854991851
$ws.Range("J24").Value = 1.048647778621615
$ws.Range("K24").Value = 1.046268937478222
$ws.Range("L24").Value = 1.057002078051427
$ws.Range("M24").Value = 1.063953984850242
$ws.Range("N24").Value = 1.050136978748839

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044457569179344
$ws.Range("D25").Value = 1.044367116971516
$ws.Range("E25").Value = 1.056041874756619
$ws.Range("F25").Value = 1.063182702050395
$ws.Range("I25").Value = 1.040878166637603
$ws.Range("J25").Value = 1.049998748801244
$ws.Range("K25").Value = 1.047004351974475
$ws.Range("L25").Value = 1.059032119487511
$ws.Range("M25").Value = 1.066151590884623
$ws.Range("N25").Value = 1.051489867461082
